$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new header cells I1 ("I0") and J1 ("IF"), matching the style already
# used by the rest of the header row (bold, bordered, centered).
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)  # xlPasteFormats

# Data for columns I (I0) and J (IF), rows 2-72.
$iValues = @(9,8,5,4,10,4,8,7,7,5,8,8,6,7,8,7,5,7,5,7,5,7,8,7,6,6,9,5,8,7,3,8,8,4,7,6,4,7,6,6,8,7,5,9,8,6,7,7,7,7,5,6,7,9,8,4,7,6,7,5,7,8,9,7,7,8,9,7,7,5,4)
$jValues = @(9,8,5,5,10,5,8,7,7,5,8,8,6,8,8,8,6,7,5,8,6,7,8,7,6,6,9,6,8,7,5,8,8,6,7,6,6,7,6,7,8,7,6,9,8,6,7,7,7,8,5,7,7,9,8,4,7,6,7,6,7,8,9,8,8,8,9,7,7,5,4)

for ($idx = 0; $idx -lt $iValues.Length; $idx++) {
    $r = $idx + 2
    $ws.Cells.Item($r, 9).Value = $iValues[$idx]
    $ws.Cells.Item($r, 10).Value = $jValues[$idx]
}
